$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$ws1 = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$ws1.Range("B9").Value = 26
$ws1.Range("C9").Value = 23
$ws1.Range("D9").Value = 36
$ws1.Range("E9").Value = 33
$ws1.Range("F9").Value = 32
$ws1.Range("G9").Value = 31

# Row 10: Total Respiratory (hazard quotient)
$ws1.Range("B10").Value = 0.31
$ws1.Range("C10").Value = 0.27
$ws1.Range("D10").Value = 0.35
$ws1.Range("E10").Value = 0.35
$ws1.Range("F10").Value = 0.35
$ws1.Range("G10").Value = 0.34

# --- Sheet "Standard Deviations" ---
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million) SD
$ws2.Range("B9").Value = 8.3
$ws2.Range("C9").Value = 7.2
$ws2.Range("D9").Value = 14
$ws2.Range("E9").Value = 12
$ws2.Range("F9").Value = 9.6
$ws2.Range("G9").Value = 9.1

# Row 10: Total Respiratory (hazard quotient) SD
$ws2.Range("B10").Value = 0.11
$ws2.Range("C10").Value = 0.094
$ws2.Range("D10").Value = 0.11
$ws2.Range("E10").Value = 0.098
$ws2.Range("F10").Value = 0.093
$ws2.Range("G10").Value = 0.089
